# 2023B Question 1 & 2: recompute the overlap-rate row (row 4, columns C:J)
# with the updated values, and drop the old "no overlap" marker row (row 5)
# that the previous (wrong) computation had tacked on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "与前一条测线的重叠率/%" values for measuring points at
# -400 .. 800 m (columns C..J). The last two survey lines now overlap fully,
# so I4/J4 become plain 0 instead of the old negative values / the "不重叠"
# text row underneath them.
$ws.Range("C4").Value = 20.9462103504215
$ws.Range("D4").Value = 18.0063169822532
$ws.Range("E4").Value = 14.8393173052909
$ws.Range("F4").Value = 11.4178381417183
$ws.Range("G4").Value = 7.70992311478901
$ws.Range("H4").Value = 3.67803151404114
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0

# Row 5 (the old "不重叠"/"不重叠" shared-string markers under I/J) is no
# longer needed now that I4/J4 are real numeric overlap rates.
$ws.Rows.Item(5).Delete()

# Leave the selection where the author's session ended up.
$ws.Range("K11").Select() | Out-Null
